$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.3920720726908886
$ws.Range("J2").Value = 0.3920720726908886
$ws.Range("M2").Value = 24.336792
$ws.Range("N2").Value = 73.01037599999999
$ws.Range("O2").Value = 0.1284781564291437
$ws.Range("P2").Value = 0.1284781564291437
$ws.Range("Q2").Value = 27.782695165128
$ws.Range("R2").Value = 250.044256486152
$ws.Range("S2").Value = 0.0503726970866786
$ws.Range("T2").Value = 0.05037269708667858

# Row 3
$ws.Range("I3").Value = 0.3920720726908886
$ws.Range("J3").Value = 0.3920720726908886
$ws.Range("O3").Value = 0.6245590297266973
$ws.Range("P3").Value = 0.6245590297266973
$ws.Range("S3").Value = 0.2448721533027565
$ws.Range("T3").Value = 0.2448721533027565

# Row 4
$ws.Range("I4").Value = 0.3920720726908886
$ws.Range("J4").Value = 0.3920720726908886
$ws.Range("M4").Value = 46.78057966666666
$ws.Range("N4").Value = 140.341739
$ws.Range("O4").Value = 0.2469628138441591
$ws.Range("P4").Value = 0.2469628138441591
$ws.Range("Q4").Value = 53.40435109635588
$ws.Range("R4").Value = 480.639159867203
$ws.Range("S4").Value = 0.09682722230145352
$ws.Range("T4").Value = 0.09682722230145351

# Row 5
$ws.Range("G5").Value = 1.770097666666667
$ws.Range("H5").Value = 5.310293000000001
$ws.Range("I5").Value = 0.6079279273091115
$ws.Range("J5").Value = 0.6079279273091115
$ws.Range("M5").Value = 24.336792
$ws.Range("N5").Value = 73.01037599999999
$ws.Range("O5").Value = 0.1284781564291437
$ws.Range("P5").Value = 0.1284781564291437
$ws.Range("Q5").Value = 43.078498733352
$ws.Range("R5").Value = 387.706488600168
$ws.Range("S5").Value = 0.07810545934246514
$ws.Range("T5").Value = 0.07810545934246513

# Row 6
$ws.Range("G6").Value = 1.770097666666667
$ws.Range("H6").Value = 5.310293000000001
$ws.Range("I6").Value = 0.6079279273091115
$ws.Range("J6").Value = 0.6079279273091115
$ws.Range("O6").Value = 0.6245590297266973
$ws.Range("P6").Value = 0.6245590297266973
$ws.Range("Q6").Value = 209.4135386027534
$ws.Range("R6").Value = 1884.721847424781
$ws.Range("S6").Value = 0.3796868764239408
$ws.Range("T6").Value = 0.3796868764239408

# Row 7
$ws.Range("G7").Value = 1.770097666666667
$ws.Range("H7").Value = 5.310293000000001
$ws.Range("I7").Value = 0.6079279273091115
$ws.Range("J7").Value = 0.6079279273091115
$ws.Range("M7").Value = 46.78057966666666
$ws.Range("N7").Value = 140.341739
$ws.Range("O7").Value = 0.2469628138441591
$ws.Range("P7").Value = 0.2469628138441591
$ws.Range("Q7").Value = 82.80619491328078
$ws.Range("R7").Value = 745.255754219527
$ws.Range("S7").Value = 0.1501355915427056
$ws.Range("T7").Value = 0.1501355915427056
